$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-5 (columns A..AH), replacing the previous
# sample readings with the new 1000-record dataset's first few rows.
$newRows = @(
    @(45132.50694444445, 10.726, 7.333, 3.404, 23.56, 17.15, 8.176, 24.228, 13.347, 5.245, 7.323, 9.309, 10.191, 2.44, 8.647, 11.655, 7.955, 2.648, 1.093, 124.223, 23.834, 7.982, 14.964, 8.049, 2.19, 13.597, 7.05, 6.629, 7.562, 9.942, 2.682, 21.628, 4.075, 9.978),
    @(45132.51388888889, 1.704, 1.01, 1.267, 4.092, 2.153, 1.222, 9.408, 2.296, 0.878, 0.635, 1.603, 1.917, 0.284, 1.504, 1.952, 1.784, 1.214, 0.341, 15.656, 4.575, 1.388, 2.617, 1.302, 0.607, 4.991, 1.226, 1.36, 1.52, 1.705, 1.136, 9.206, 0.514, 1.743),
    @(45132.52083333334, 12.808, 9.473, 1.183, 28.182, 22.464, 10.003, 35.281, 15.677, 6.888, 9.892, 11.292, 12.112, 3.092, 10.151, 14.244, 8.795, 0.906, 0.544, 147.148, 28.216, 9.37, 18.686, 9.872, 1.584, 17.893, 8.276, 7.454, 8.733, 11.878, 0.722, 31.776, 5.148, 11.714),
    @(45132.52777777778, 7.08, 5.22, 0.79, 15.66, 12.31, 5.51, 23.97, 8.7, 3.85, 5.36, 6.27, 6.77, 1.68, 5.64, 7.91, 4.97, 0.66, 0.32, 78.52, 15.88, 5.21, 10.46, 5.47, 0.92, 11.73, 4.6, 4.19, 4.89, 6.6, 0.53, 22.03, 2.82, 6.51)
)

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowValues = $newRows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowValues[$c]
    }
}

# Row 6 (an extra sample row) is no longer part of the dataset; remove it so
# the sheet's used range shrinks back to A1:AH5.
$ws.Rows.Item(6).Delete()

# A handful of columns were re-measured for "custom accuracy" and now use a
# slightly different width (in Excel "characters"); ColumnWidth must be set
# to targetWidth - 5/6 because Excel adds a fixed ~0.8333 character padding
# when translating ColumnWidth into the stored column width.
$widthAdjust = 5.0 / 6.0
$colWidths = @{
    3  = 7;
    11 = 7;
    13 = 8;
    22 = 7;
    24 = 7;
    34 = 8;
}

foreach ($colIndex in $colWidths.Keys) {
    $ws.Columns.Item($colIndex).ColumnWidth = $colWidths[$colIndex] - $widthAdjust
}
